$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.674.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.645.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.662.02"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.87%  "
$ws.Range("E10").Value = "  +3.24%  "
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.337"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.119.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.598.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.640.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "344.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.415"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0751"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("E36").Value = "  +2.22%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.838"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.80%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.824"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "291.53"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.73%  "
$ws.Range("E41").Value = "  +1.59%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("E44").Value = "  +4.19%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0950"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.968.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.89%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.16%  "
